# Renumber the path-index column (A) in each of the four LSMC path tables
# from a 1-based index (1..8) to a 0-based index (0..7), matching the
# "renumbered the paths from 0" fix described in the commit message.
#
# Affected blocks on sheet "Original":
#   A37:A44  (t3 payoff table)
#   A60:A67  (t2 regression table)   -- L60:L64 is a FILTER(A60:A67,...)
#            spill/array formula, so its cached values shift automatically
#            on recalculation once column A is corrected.
#   A75:A82  (t1 regression table)   -- L75:L79 is the analogous FILTER spill.
#   A91:A98  (final payoff table)
#
# Also updates the sheet selection to match the author's final cursor
# position/selection after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Original")
if (-not $ws) { $ws = $wb.ActiveSheet }

function Set-PathIndex($startRow, $endRow) {
    $newValue = 0
    for ($row = $startRow; $row -le $endRow; $row++) {
        $ws.Cells.Item($row, 1).Value = $newValue
        $newValue = $newValue + 1
    }
}

Set-PathIndex 37 44
Set-PathIndex 60 67
Set-PathIndex 75 82
Set-PathIndex 91 98

# Match the author's final selection/cursor (rows 91:98, active cell A91)
# and scroll the view down so row 68 is the top visible row (best-effort;
# some hosts don't expose/persist window-scroll state).
$ws.Range("A91:A98").Select()
try {
    $excel.ActiveWindow.ScrollRow = 68
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
